# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B:G, rows 2-10 (TB, d2S, K, IP, Win, sum)
$data = @{
    2  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 8.656069925401464)
    3  = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 3.755628166162433)
    4  = @(0.6545652718822623, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 1, 1.950081467445961)
    5  = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 24.14949828602258)
    6  = @(3.272327238179451, 1.626987699542094, 189.6080260415259, 13.86384647080068, 0, 208.3711874500482)
    7  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 13.86384647080068, 1, 21.98653043760045)
    8  = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 24.14949828602258)
    9  = @(0.2881169905109251, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 0, 1.583633186074624)
    10 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 8.656069925401464)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B: TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C: d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D: K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E: IP
    $ws.Cells.Item($row, 6).Value = $vals[4]  # F: Win
    $ws.Cells.Item($row, 7).Value = $vals[5]  # G: sum
}
